# Add two new columns, "I0" (column I) and "IF" (column J), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: copy H1's formatting (bold font, border, center/top align)
# onto I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for the new columns (row -> [I value, J value]).
$values = @{
    2 = @(5, 7)
    3 = @(4, 7)
    4 = @(4, 6)
    5 = @(6, 8)
    6 = @(8, 9)
    7 = @(1, 3)
    8 = @(7, 8)
    9 = @(5, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
